# Revert capacity charts to show kilowatts on the y-axis.
# 1. Convert the Watts values in column E (rows 11-26) to Kilowatts (divide by 1000).
# 2. Update the number format used by those cells to show one decimal place.
# 3. Update the chart's value-axis title and number format to reflect kilowatts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1 & 2: convert E11:E26 from Watts to Kilowatts -------------------------
for ($row = 11; $row -le 26; $row++) {
    $cell = $ws.Cells.Item($row, 5)   # column E
    [double]$watts = $cell.Value2
    $cell.Value = $watts / 1000
}

# Update the custom number format (numFmtId 164) so it shows one decimal
# place, e.g. "#,##0" -> "#,##0.0". This format is shared by every numeric
# data cell on the sheet (columns B-G), so update the whole body.
$ws.Range("B2:G26").NumberFormat = "#,##0.0"

# --- 3: update the chart ----------------------------------------------------
$chartObj = $ws.ChartObjects(1)
$chart = $chartObj.Chart
$valueAxis = $chart.Axes(2)          # 2 = xlValue

$valueAxis.AxisTitle.Text = "Kilowatts (kW)"
$valueAxis.TickLabels.NumberFormat = "#,##0"
